# Apply the "riska.xlsx" update: refreshed repayment figures for the
# 2025-09-23 cycle, re-exported/re-uploaded (sheet tab bumped from (2) to (4)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab to match the re-uploaded export ---
$ws.Name = "repayment_20250923_20250923 (4)"

# --- Updated repayment figures per collector (row = collector) ---
# Row 2: Annisa Putri Restu - Talk_time refreshed
$ws.Range("H2").Value = 1.149

# Row 3: Azizah Rahmawati - Talk_time refreshed
$ws.Range("H3").Value = 227

# Row 4: Erlangga Hutama
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "3,143,300.00"
$ws.Range("G4").Value = "2.11"
$ws.Range("H4").Value = 51

# Row 5: Aldi Taufik
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = "4,457,357.00"
$ws.Range("G5").Value = "2.68"
$ws.Range("H5").Value = 637

# Row 6: Yandi Nugraha
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "942,556.00"
$ws.Range("G6").Value = "0.71"
$ws.Range("H6").Value = 54

# Row 7: Ridhoi Berkat Zebua - Talk_time refreshed
$ws.Range("H7").Value = 1.354

# Row 8: Riska Nurlita - Talk_time refreshed
$ws.Range("H8").Value = 493

# Row 9: Debora Retima Sihombing
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = "479,619.00"
$ws.Range("G9").Value = "0.32"
$ws.Range("H9").Value = 664

# Row 10: Erick Ervan Dewanggga
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "452,511.00"
$ws.Range("G10").Value = "0.29"
$ws.Range("H10").Value = 235

# Row 11: Romli
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "747,151.00"
$ws.Range("G11").Value = "0.53"
$ws.Range("H11").Value = 816

# Row 12: Fadilah Damayanti
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = "450,000.00"
$ws.Range("G12").Value = "0.33"
$ws.Range("H12").Value = 371
$ws.Range("K12").Value = "3.34"

# Row 13: Nur Halim
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "344,060.00"
$ws.Range("G13").Value = "0.25"
$ws.Range("H13").Value = 109

# Row 14: Adistira Winditya P
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "726,691.00"
$ws.Range("G14").Value = "0.46"
$ws.Range("H14").Value = 44

# Row 15: Axl Wicaksono - Talk_time refreshed
$ws.Range("H15").Value = 97

# Row 16: Sucika Wardani
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "527,055.00"
$ws.Range("G16").Value = "0.32"
$ws.Range("H16").Value = 37

# Row 17: Wasti Feronika Sihombing - Talk_time refreshed
$ws.Range("H17").Value = 272

# Row 18: Nuraini - Talk_time refreshed
$ws.Range("H18").Value = 259

# --- Resize columns to fit the refreshed content ---
$ws.Columns("A:L").AutoFit()
